$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Locate the split point right after "Available" (before " at: ")
#    in the "Available at: <link>" paragraph, and drop a fresh
#    "_GoBack" bookmark there. Word only ever keeps a single
#    "_GoBack" bookmark, so adding this one automatically removes
#    the old one that sat in the title paragraph (around "VI").
# ------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("Available", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $findRange.End

$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 2) Now swap the word "Available" for "Try them" (the trailing
#    " at: " text and the hyperlink that follows are untouched).
#    Doing this after the bookmark split lets the now-isolated run
#    recompute whether it still needs xml:space="preserve".
# ------------------------------------------------------------------
$replaceRange = $d.Content
$replaceRange.Find.Execute("Available", $false, $false, $false, $false, $false, $true, 1, $false, "Try them", 2)
